# 自动更新Excel文件 - 2026-02-26 23:22:21
# Recompute the "剩余" (remaining days) and "开始时间" (start date) columns
# as of the new reference date 2026-02-27, for every data row.
#   剩余(E) = 总天(D) - (今天 - 开始时间(F))
# When the remaining days drop to zero or below, the cycle restarts:
#   开始时间(F) = 今天, 剩余(E) = 总天(D)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = [datetime]::ParseExact("20260227", "yyyyMMdd", $null)
$todayNum = [int]$today.ToString("yyyyMMdd")
$todayOA = $today.ToOADate()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value2
    $f = $fCell.Value2

    if ($d -eq $null -or $f -eq $null) {
        continue
    }

    $fstr = [string]([int64]$f)
    if ($fstr.Length -ne 8) {
        continue
    }

    $parsedOk = $true
    try {
        $fdate = [datetime]::ParseExact($fstr, "yyyyMMdd", $null)
    } catch {
        $parsedOk = $false
    }
    if (-not $parsedOk) {
        continue
    }

    $fOA = $fdate.ToOADate()
    $elapsed = [int]($todayOA - $fOA)
    $newE = [int]$d - $elapsed

    if ($newE -le 0) {
        $newE = [int]$d
        $fCell.Value2 = $todayNum
    }
    $eCell.Value2 = $newE
}
